$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.24677300453186
$ws.Range("B1").Value = 2.648993015289307
$ws.Range("C1").Value = 8.482213973999023
$ws.Range("D1").Value = 2.113804817199707
$ws.Range("E1").Value = 1.141023635864258
